$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric but must remain literal text
# (matches source Price formatting, e.g. trailing zeros / 3-decimal coin prices).
# Force a Text number format on just these cells before assigning, so Excel
# stores them as strings instead of re-parsing as numbers.
$textCells = @("D4", "D5", "D7", "D8", "D9", "D14", "D15", "D16", "D18", "D20", "D23", "D24", "D25", "D26", "D27", "D35", "D37", "D40", "D42", "D43", "D46", "D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.345.75"
$ws.Range("E2").Value = "  -1.85%  "

$ws.Range("D3").Value = "1.655.79"
$ws.Range("E3").Value = "  -0.70%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "213.26"
$ws.Range("E5").Value = "  -0.66%  "

$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("D8").Value = "23.70"
$ws.Range("E8").Value = "  +0.57%  "

$ws.Range("D9").Value = "0.261"
$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("E10").Value = "  -1.04%  "

$ws.Range("E11").Value = "  -0.34%  "

$ws.Range("D12").Value = "1.888.87"
$ws.Range("E12").Value = "  -0.86%  "

$ws.Range("D13").Value = "1.651.66"
$ws.Range("E13").Value = "  -1.03%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "4.08"
$ws.Range("E14").Value = "  -1.71%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.572"
$ws.Range("E15").Value = "  +3.65%  "

$ws.Range("D16").Value = "65.80"
$ws.Range("E16").Value = "  -0.31%  "

$ws.Range("D17").Value = "27.344.41"
$ws.Range("E17").Value = "  -1.85%  "

$ws.Range("D18").Value = "232.46"
$ws.Range("E18").Value = "  -7.37%  "

$ws.Range("E19").Value = "  -0.69%  "

$ws.Range("D20").Value = "7.48"
$ws.Range("E20").Value = "  -0.98%  "

$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("E22").Value = "  -2.17%  "

$ws.Range("D23").Value = "9.38"
$ws.Range("E23").Value = "  +0.38%  "

$ws.Range("D24").Value = "2.02"
$ws.Range("E24").Value = "  -1.95%  "

$ws.Range("D25").Value = "147.15"
$ws.Range("E25").Value = "  +0.28%  "

$ws.Range("D26").Value = "7.15"
$ws.Range("E26").Value = "  -1.00%  "

$ws.Range("D27").Value = "15.89"
$ws.Range("E27").Value = "  -2.42%  "

$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("E29").Value = "  -0.61%  "

$ws.Range("E30").Value = "  -1.06%  "

$ws.Range("E31").Value = "  -3.73%  "

$ws.Range("E32").Value = "  -1.31%  "

$ws.Range("D33").Value = "1.450.75"
$ws.Range("E33").Value = "  +1.48%  "

$ws.Range("E34").Value = "  -0.40%  "

$ws.Range("D35").Value = "1.56"
$ws.Range("E35").Value = "  -0.21%  "

$ws.Range("E36").Value = "  -0.72%  "

$ws.Range("D37").Value = "0.910"
$ws.Range("E37").Value = "  -2.14%  "

$ws.Range("E38").Value = "  -1.79%  "

$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("D40").Value = "1.04"
$ws.Range("E40").Value = "  +0.69%  "

$ws.Range("E41").Value = "  -0.17%  "

$ws.Range("D42").Value = "5.52"
$ws.Range("E42").Value = "  +2.45%  "

$ws.Range("D43").Value = "65.24"
$ws.Range("E43").Value = "  -6.18%  "

$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("D45").Value = "1.796.52"
$ws.Range("E45").Value = "  -0.86%  "

$ws.Range("D46").Value = "0.786"
$ws.Range("E46").Value = "  -0.59%  "

$ws.Range("E47").Value = "  -0.42%  "

$ws.Range("D48").Value = "88.20"
$ws.Range("E48").Value = "  -0.83%  "

$ws.Range("E49").Value = "  +0.56%  "

$ws.Range("E50").Value = "  -0.35%  "

$ws.Range("E51").Value = "  -1.19%  "
